$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 180.8
$ws.Range("B8").Value = 1.25
$ws.Range("B10").Value = 10
$ws.Range("B13").Value = "25,66,86"
$ws.Range("B15").Value = 0.25
$ws.Range("B27").Value = 0.33

$ws.Range("C18").Select()
